# Update Leve profit-tracking sheets with refreshed market-board pricing data
# (scheduled runner sync across all job sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 4000
$ws.Range("J44").Value = 4000
$ws.Range("L44").Value = 4000
$ws.Range("N44").Value = -4924
$ws.Range("H51").Value = 2816.0952
$ws.Range("I51").Value = 2276.4443
$ws.Range("K51").Value = 2276.4443
$ws.Range("M51").Value = -1792.4443
$ws.Range("H107").Value = 1872.1052
$ws.Range("I107").Value = 1759.4667
$ws.Range("J107").Value = 2294.5
$ws.Range("K107").Value = 1759.4667
$ws.Range("L107").Value = 2294.5
$ws.Range("M107").Value = 160.5333000000001
$ws.Range("N107").Value = -6134.5
$ws.Range("H116").Value = 3382.762
$ws.Range("I116").Value = 3056
$ws.Range("J116").Value = 3679.818
$ws.Range("K116").Value = 3056
$ws.Range("L116").Value = 3679.818
$ws.Range("M116").Value = 386
$ws.Range("N116").Value = -10563.818
$ws.Range("H137").Value = 1058.2941
$ws.Range("I137").Value = 886.8333
$ws.Range("J137").Value = 1469.8
$ws.Range("K137").Value = 2660.4999
$ws.Range("L137").Value = 4409.4
$ws.Range("M137").Value = -110.4998999999998
$ws.Range("N137").Value = -9509.4
$ws.Range("H138").Value = 4741.7085
$ws.Range("I138").Value = 4779.25
$ws.Range("J138").Value = 4722.9375
$ws.Range("K138").Value = 14337.75
$ws.Range("L138").Value = 14168.8125
$ws.Range("M138").Value = -9197.75
$ws.Range("N138").Value = -24448.8125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 395716.38
$ws.Range("I32").Value = 3380.1147
$ws.Range("J32").Value = 2390092.2
$ws.Range("K32").Value = 3380.1147
$ws.Range("L32").Value = 2390092.2
$ws.Range("M32").Value = -3093.1147
$ws.Range("N32").Value = -2390666.2
$ws.Range("H61").Value = 2400
$ws.Range("I61").Value = 2750
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 2750
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -2538
$ws.Range("N61").Value = -1424
$ws.Range("H74").Value = 1989
$ws.Range("J74").Value = 2275
$ws.Range("L74").Value = 2275
$ws.Range("N74").Value = -4023
$ws.Range("H77").Value = 1989
$ws.Range("J77").Value = 2275
$ws.Range("L77").Value = 11375
$ws.Range("N77").Value = -20111
$ws.Range("H110").Value = 2021.1818
$ws.Range("I110").Value = 1849.3334
$ws.Range("J110").Value = 2794.5
$ws.Range("K110").Value = 1849.3334
$ws.Range("L110").Value = 2794.5
$ws.Range("M110").Value = 195.6666
$ws.Range("N110").Value = -6884.5
$ws.Range("H122").Value = 28835.652
$ws.Range("I122").Value = 32781.65
$ws.Range("J122").Value = 2529
$ws.Range("K122").Value = 98344.95000000001
$ws.Range("L122").Value = 7587
$ws.Range("M122").Value = -95894.95000000001
$ws.Range("N122").Value = -12487
$ws.Range("H132").Value = 1942.2142
$ws.Range("I132").Value = 1188.5927
$ws.Range("J132").Value = 3298.7334
$ws.Range("K132").Value = 3565.7781
$ws.Range("L132").Value = 9896.200199999999
$ws.Range("M132").Value = -1035.7781
$ws.Range("N132").Value = -14956.2002
$ws.Range("H134").Value = 42900
$ws.Range("J134").Value = 42900
$ws.Range("L134").Value = 42900
$ws.Range("N134").Value = -53040
$ws.Range("H135").Value = 69429
$ws.Range("J135").Value = 69429
$ws.Range("L135").Value = 69429
$ws.Range("N135").Value = -79569
$ws.Range("H136").Value = 2400
$ws.Range("I136").Value = 2750
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 8250
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -5700
$ws.Range("N136").Value = -8100
$ws.Range("H137").Value = 86667
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19098.125
$ws.Range("I134").Value = 1330.8334
$ws.Range("J134").Value = 72400
$ws.Range("K134").Value = 3992.5002
$ws.Range("L134").Value = 217200
$ws.Range("M134").Value = -1457.5002
$ws.Range("N134").Value = -222270
$ws.Range("H138").Value = 19250

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1988.2745
$ws.Range("I31").Value = 1556.3429
$ws.Range("J31").Value = 2933.125
$ws.Range("K31").Value = 1556.3429
$ws.Range("L31").Value = 2933.125
$ws.Range("M31").Value = -1261.3429
$ws.Range("N31").Value = -3523.125
$ws.Range("H34").Value = 1988.2745
$ws.Range("I34").Value = 1556.3429
$ws.Range("J34").Value = 2933.125
$ws.Range("K34").Value = 1556.3429
$ws.Range("L34").Value = 2933.125
$ws.Range("M34").Value = -1354.3429
$ws.Range("N34").Value = -3337.125
$ws.Range("H58").Value = 1224.6818
$ws.Range("I58").Value = 991.7368
$ws.Range("J58").Value = 2700
$ws.Range("K58").Value = 991.7368
$ws.Range("L58").Value = 2700
$ws.Range("M58").Value = -788.7368
$ws.Range("N58").Value = -3106
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H136").Value = 1224.6818
$ws.Range("I136").Value = 991.7368
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 2975.2104
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -425.2103999999999
$ws.Range("N136").Value = -13200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 33713.84
$ws.Range("I121").Value = 4475.7144
$ws.Range("J121").Value = 42241.625
$ws.Range("K121").Value = 13427.1432
$ws.Range("L121").Value = 126724.875
$ws.Range("M121").Value = -12117.1432
$ws.Range("N121").Value = -129344.875
$ws.Range("H131").Value = 5883217
$ws.Range("J131").Value = 6494338
$ws.Range("L131").Value = 19483014
$ws.Range("N131").Value = -19493094
$ws.Range("H132").Value = 1184.8438
$ws.Range("I132").Value = 1113.0555
$ws.Range("J132").Value = 1277.1428
$ws.Range("K132").Value = 10017.4995
$ws.Range("L132").Value = 11494.2852
$ws.Range("M132").Value = -7487.4995
$ws.Range("N132").Value = -16554.2852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4185
$ws.Range("I18").Value = 4185
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4185
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3892
$ws.Range("N18").Value = ""
$ws.Range("H21").Value = 1795.7273
$ws.Range("I21").Value = 502
$ws.Range("K21").Value = 502
$ws.Range("M21").Value = -329
$ws.Range("H30").Value = 1795.7273
$ws.Range("I30").Value = 502
$ws.Range("K30").Value = 502
$ws.Range("M30").Value = -397
$ws.Range("H122").Value = 4238.909
$ws.Range("I122").Value = 4958.5
$ws.Range("J122").Value = 2320
$ws.Range("K122").Value = 14875.5
$ws.Range("L122").Value = 6960
$ws.Range("M122").Value = -12425.5
$ws.Range("N122").Value = -11860
$ws.Range("H132").Value = 4529.4883
$ws.Range("I132").Value = 4655.029
$ws.Range("J132").Value = 3980.25
$ws.Range("K132").Value = 13965.087
$ws.Range("L132").Value = 11940.75
$ws.Range("M132").Value = -11435.087
$ws.Range("N132").Value = -17000.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2807.2
$ws.Range("J122").Value = 2862.3333
$ws.Range("L122").Value = 8586.999899999999
$ws.Range("N122").Value = -13486.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 175001780
$ws.Range("I2").Value = 350000930
$ws.Range("J2").Value = 2626.6667
$ws.Range("K2").Value = 350000930
$ws.Range("L2").Value = 2626.6667
$ws.Range("M2").Value = -350000818
$ws.Range("N2").Value = -2850.6667
$ws.Range("H15").Value = 7606
$ws.Range("I15").Value = 7206
$ws.Range("J15").Value = 8006
$ws.Range("K15").Value = 7206
$ws.Range("L15").Value = 8006
$ws.Range("M15").Value = -6918
$ws.Range("N15").Value = -8582
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1710
$ws.Range("N29").Value = ""
$ws.Range("H122").Value = 1785.4286
$ws.Range("I122").Value = 1687.8823
$ws.Range("K122").Value = 5063.6469
$ws.Range("M122").Value = -2613.6469
$ws.Range("H136").Value = 872.38464
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

Write-Output "Applied profit updates across all sheets"